# Regenerate merged AHB files
#
# 1) Rename the "_old"/"_new" suffixed header cells in row 1 to
#    "_FV2310"/"_FV2404" respectively (the "diff" header is unchanged).
# 2) Turn the used range A1:U68 into a native Excel table ("Table1")
#    with an AutoFilter, whose column names are taken from the
#    (already renamed) header row.
# 3) Freeze the header row (split after row 1, active pane bottom-left).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Create the table over the full used range, using the header row
# that was just renamed, and keep the default (no) style so no extra
# conditional-format / style definitions are introduced.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U68"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze the top row (header).
$ws.Range("A2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
